$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the last existing data row (196) as a template so the new rows pick up
# the same per-row cell layout (including the otherwise-empty K and M cells)
# as the rest of the table.
$ws.Range("A196:M196").Copy($ws.Range("A197:M197"))
$ws.Range("A196:M196").Copy($ws.Range("A198:M198"))
$ws.Range("A196:M196").Copy($ws.Range("A199:M199"))

# Row 197 (record #196): 9U423 Chisinau (KIV) Air Moldova A319 (ER-AXL)
$ws.Cells.Item(197, 1).Value = 196
$ws.Cells.Item(197, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(197, 3).Value = "8:00 AM"
$ws.Cells.Item(197, 4).Value = "9U423"
$ws.Cells.Item(197, 5).Value = "Chisinau"
$ws.Cells.Item(197, 6).Value = "(KIV)"
$ws.Cells.Item(197, 7).Value = "Air Moldova "
$ws.Cells.Item(197, 8).Value = "A319"
$ws.Cells.Item(197, 9).Value = "(ER-AXL)"
$ws.Cells.Item(197, 10).Value = "7:30 AM"
$ws.Cells.Item(197, 12).Value = "0 hours, -30 minutes"

# Row 198 (record #197): FR4669 Stockholm (ARN) Ryanair B738 (9H-QDK)
$ws.Cells.Item(198, 1).Value = 197
$ws.Cells.Item(198, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(198, 3).Value = "9:10 AM"
$ws.Cells.Item(198, 4).Value = "FR4669"
$ws.Cells.Item(198, 5).Value = "Stockholm"
$ws.Cells.Item(198, 6).Value = "(ARN)"
$ws.Cells.Item(198, 7).Value = "Ryanair "
$ws.Cells.Item(198, 8).Value = "B738"
$ws.Cells.Item(198, 9).Value = "(9H-QDK)"
$ws.Cells.Item(198, 10).Value = "8:55 AM"
$ws.Cells.Item(198, 12).Value = "0 hours, -15 minutes"

# Row 199 (record #198): FR8083 Birmingham (BHX) Ryanair B738 (EI-EXD)
$ws.Cells.Item(199, 1).Value = 198
$ws.Cells.Item(199, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(199, 3).Value = "9:50 AM"
$ws.Cells.Item(199, 4).Value = "FR8083"
$ws.Cells.Item(199, 5).Value = "Birmingham"
$ws.Cells.Item(199, 6).Value = "(BHX)"
$ws.Cells.Item(199, 7).Value = "Ryanair "
$ws.Cells.Item(199, 8).Value = "B738"
$ws.Cells.Item(199, 9).Value = "(EI-EXD)"
$ws.Cells.Item(199, 10).Value = "9:47 AM"
$ws.Cells.Item(199, 12).Value = "0 hours, -3 minutes"
